$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column D width (closest achievable value to 18.125 given engine quantization)
$ws.Columns.Item(4).ColumnWidth = 17.4

# Update character roster: ids, names, and completion status (v=voiced, c=complete, x=broken)
$ws.Cells.Item(1,1).Value = 1001
$ws.Cells.Item(1,2).Value = "Hiyori"
$ws.Cells.Item(1,3).Value = "v"
$ws.Cells.Item(2,1).Value = 1002
$ws.Cells.Item(2,2).Value = "Yui"
$ws.Cells.Item(2,3).Value = "v"
$ws.Cells.Item(3,1).Value = 1003
$ws.Cells.Item(3,2).Value = "Rei"
$ws.Cells.Item(3,3).Value = "v"
$ws.Cells.Item(4,1).Value = 1004
$ws.Cells.Item(4,2).Value = "Misogi"
$ws.Cells.Item(4,3).Value = "v"
$ws.Cells.Item(5,1).Value = 1006
$ws.Cells.Item(5,2).Value = "Akari"
$ws.Cells.Item(5,3).Value = "v"
$ws.Cells.Item(6,1).Value = 1007
$ws.Cells.Item(6,2).Value = "Miyako"
$ws.Cells.Item(6,3).Value = "v"
$ws.Cells.Item(7,1).Value = 1008
$ws.Cells.Item(7,2).Value = "Yuki"
$ws.Cells.Item(7,3).Value = "v"
$ws.Cells.Item(8,1).Value = 1009
$ws.Cells.Item(8,2).Value = "Anna"
$ws.Cells.Item(8,3).Value = "v"
$ws.Cells.Item(9,1).Value = 1010
$ws.Cells.Item(9,2).Value = "Maho"
$ws.Cells.Item(9,3).Value = "c"
$ws.Cells.Item(10,1).Value = 1011
$ws.Cells.Item(10,2).Value = "Rino"
$ws.Cells.Item(10,3).Value = "c"
$ws.Cells.Item(11,1).Value = 1012
$ws.Cells.Item(11,2).Value = "Hatsune"
$ws.Cells.Item(11,3).Value = "c"
$ws.Cells.Item(12,1).Value = 1016
$ws.Cells.Item(12,2).Value = "Suzuna"
$ws.Cells.Item(12,3).Value = "c"
$ws.Cells.Item(13,1).Value = 1017
$ws.Cells.Item(13,2).Value = "Kaori"
$ws.Cells.Item(13,3).Value = "c"
$ws.Cells.Item(14,1).Value = 1018
$ws.Cells.Item(14,2).Value = "Io"
$ws.Cells.Item(14,3).Value = "c"
$ws.Cells.Item(15,1).Value = 1020
$ws.Cells.Item(15,2).Value = "Mimi"
$ws.Cells.Item(15,3).Value = "c"
$ws.Cells.Item(16,1).Value = 1021
$ws.Cells.Item(16,2).Value = "Kurumi"
$ws.Cells.Item(16,3).Value = "c"
$ws.Cells.Item(17,1).Value = 1022
$ws.Cells.Item(17,2).Value = "Yori"
$ws.Cells.Item(17,3).Value = "c"
$ws.Cells.Item(18,1).Value = 1025
$ws.Cells.Item(18,2).Value = "Suzume"
$ws.Cells.Item(18,3).Value = "c"
$ws.Cells.Item(19,1).Value = 1027
$ws.Cells.Item(19,2).Value = "Eriko"
$ws.Cells.Item(19,3).Value = "c"
$ws.Cells.Item(20,1).Value = 1028
$ws.Cells.Item(20,2).Value = "Saren"
$ws.Cells.Item(20,3).Value = "c"
$ws.Cells.Item(21,1).Value = 1029
$ws.Cells.Item(21,2).Value = "Nozomi"
$ws.Cells.Item(21,3).Value = "c"
$ws.Cells.Item(22,1).Value = 1030
$ws.Cells.Item(22,2).Value = "Ninon"
$ws.Cells.Item(22,3).Value = "c"
$ws.Cells.Item(23,1).Value = 1031
$ws.Cells.Item(23,2).Value = "Shinobu"
$ws.Cells.Item(23,3).Value = "c"
$ws.Cells.Item(24,1).Value = 1032
$ws.Cells.Item(24,2).Value = "Akino"
$ws.Cells.Item(24,3).Value = "c"
$ws.Cells.Item(25,1).Value = 1033
$ws.Cells.Item(25,2).Value = "Mahiru"
$ws.Cells.Item(25,3).Value = "c"
$ws.Cells.Item(26,1).Value = 1034
$ws.Cells.Item(26,2).Value = "Yukari"
$ws.Cells.Item(26,3).Value = "c"
$ws.Cells.Item(27,1).Value = 1038
$ws.Cells.Item(27,2).Value = "Shiori"
$ws.Cells.Item(27,3).Value = "c"
$ws.Cells.Item(28,1).Value = 1040
$ws.Cells.Item(28,2).Value = "Aoi"
$ws.Cells.Item(28,3).Value = "c"
$ws.Cells.Item(29,1).Value = 1042
$ws.Cells.Item(29,2).Value = "Chika"
$ws.Cells.Item(29,3).Value = "c"
$ws.Cells.Item(30,1).Value = 1043
$ws.Cells.Item(30,2).Value = "Makoto"
$ws.Cells.Item(30,3).Value = "c"
$ws.Cells.Item(31,1).Value = 1045
$ws.Cells.Item(31,2).Value = "Kuuka"
$ws.Cells.Item(31,3).Value = "c"
$ws.Cells.Item(32,1).Value = 1046
$ws.Cells.Item(32,2).Value = "Tamaki"
$ws.Cells.Item(32,3).Value = "c"
$ws.Cells.Item(33,1).Value = 1048
$ws.Cells.Item(33,2).Value = "Mifuyu"
$ws.Cells.Item(33,3).Value = "c"
$ws.Cells.Item(33,4).Value = "skill0 problem"
$ws.Cells.Item(34,1).Value = 1049
$ws.Cells.Item(34,2).Value = "Shizuru"
$ws.Cells.Item(34,3).Value = "c"
$ws.Cells.Item(35,1).Value = 1050
$ws.Cells.Item(35,2).Value = "Misaki"
$ws.Cells.Item(35,3).Value = "c"
$ws.Cells.Item(36,1).Value = 1051
$ws.Cells.Item(36,2).Value = "Mitsuki"
$ws.Cells.Item(36,3).Value = "c"
$ws.Cells.Item(37,1).Value = 1052
$ws.Cells.Item(37,2).Value = "Rima"
$ws.Cells.Item(37,3).Value = "x"
$ws.Cells.Item(37,4).Value = "skill1 broken"
$ws.Cells.Item(38,1).Value = 1053
$ws.Cells.Item(38,2).Value = "Monika"
$ws.Cells.Item(38,3).Value = "c"
$ws.Cells.Item(38,4).Value = "warn: skill0 超大"
$ws.Cells.Item(39,1).Value = 1057
$ws.Cells.Item(39,2).Value = "Djeeta"
$ws.Cells.Item(39,3).Value = "c"
$ws.Cells.Item(40,1).Value = 1058
$ws.Cells.Item(40,2).Value = "Pecorine"
$ws.Cells.Item(40,3).Value = "v"
$ws.Cells.Item(41,1).Value = 1059
$ws.Cells.Item(41,2).Value = "Kokoro"
$ws.Cells.Item(41,3).Value = "v"
$ws.Cells.Item(42,1).Value = 1060
$ws.Cells.Item(42,2).Value = "Kyaru"
$ws.Cells.Item(42,3).Value = "v"
$ws.Cells.Item(43,1).Value = 1063
$ws.Cells.Item(43,2).Value = "Arisa"
$ws.Cells.Item(43,3).Value = "v"
$ws.Cells.Item(44,3).Value = "63 56"

# Restore selection to the last-edited cell
$ws.Range("D40").Select()